$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Reword the first comment line.
$ws.Range("A1").Value = "Put the path to images under the corresponding title."

# 2) Insert a brand-new row at position 2 to hold a new explanatory sentence;
#    this pushes the former row 2 (root-folder comment), row 3 (path) and
#    row 5 (Pol0_90/Pol45_135 header) down by one row each.
$ws.Rows.Item(2).Insert()

# Merge the new row's A:F range like the other comment rows, then touch an
# inert formatting property on the B:F cells so the empty cells are written
# out (matching the structure of the sibling comment rows) before filling in
# the text in A2.
$ws.Range("A2:F2").Merge()
$ws.Range("B2:F2").Borders.Item(1).LineStyle = -4142
$ws.Range("A2").Value = "The files in each row must correspond to different polarizations of same sample."
